# Updated cryptos list on Mon Nov 25 08:39:01 UTC 2024 with GitHub Actions
# Refresh price / 1h-volume figures scraped from coinranking.com, and fix a
# handful of rows whose coin pairs had swapped order versus the previous run.
#
# Cells are plain text (t="inlineStr") in the workbook -- numeric-looking
# values are written with a leading apostrophe so Excel stores them as text
# instead of silently coercing to Number (which would also rewrite "665.10"
# as 665.1, drop the thousand-dot grouping in e.g. "98.228.59", or flip tiny
# decimals like "0.0000202" into scientific notation). The Style reset right
# after each write clears the quote-prefix marker so the cell keeps the
# workbook default (unstyled) formatting, matching the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'98.228.59"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.01%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.406.60"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.27%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.03%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'254.44"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.43%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'665.10"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -2.69%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'1.52"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +4.73%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.433"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -0.30%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'1.05"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.84%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.00%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'3.401.50"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.34%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'45.15"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +7.62%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.210"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -2.99%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'98.026.61"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.05%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'6.16"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -3.22%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  -3.01%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'4.051.83"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.21%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'9.13"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.31%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'3.413.77"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.08%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'18.34"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +4.06%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'0.547"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -5.81%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'11.38"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +2.46%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'512.40"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.39%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'3.42"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.73%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.0000202"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -1.81%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'6.81"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +3.39%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'97.70"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -3.28%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'12.44"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -3.01%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'3.588.77"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.40%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'12.11"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +3.75%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = "'PancakeSwap"
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = "'2.89"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +8.68%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = "'Hedera"
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'0.144"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -4.29%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.08%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.188"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -3.93%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.01%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  -1.81%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'29.18"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -2.48%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'7.95"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.79%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'1.49"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -3.93%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = "'Kaspa"
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = "'0.154"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.23%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = "'Bittensor"
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = "'525.69"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -1.69%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.01%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'WhiteBITCoin"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'24.42"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -1.19%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'ARBITRUM"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'0.863"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -2.08%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'1.74"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -1.27%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.0427"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -2.26%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  -2.04%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = "'Filecoin"
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'5.65"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.76%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'OKB"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'56.36"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +0.66%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -4.61%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'Stacks"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'2.23"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +4.89%  "
$ws.Range('E51').Style = 'Normal'
